# Patch class-diagram constant names (camelCase field names -> UPPER_SNAKE
# constant names) in the two UML "attribute" text boxes of the class
# diagram, mirroring a manual rename done in PowerPoint. Replacing the
# text changes each box's rendered width slightly (spAutoFit), so the
# handful of neighbouring shapes that PowerPoint nudges as a result are
# corrected too.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Replace the leading "-<camelCase>: " run-cluster of a paragraph with a
# single run containing "-<UPPER_SNAKE>: ", leaving the trailing type
# run (e.g. "int") untouched -- matches how PowerPoint merges a manual
# retype of a text selection into one run using the selection's starting
# format.
function Set-ParaPrefix {
    param($TextRange, [int]$ParaIndex, [string]$NewPrefix, [int]$OldPrefixLen)

    $para = $TextRange.Paragraphs($ParaIndex, 1)
    $chars = $TextRange.Characters($para.Start, $OldPrefixLen)
    $chars.Text = $NewPrefix
}

# --- Shape "ZoneTexte 7" (Animal class attributes) ---
$animalAttrs = $s.Shapes.Item(2)
if ($animalAttrs.Name -ne "ZoneTexte 7") { throw "unexpected shape 2: $($animalAttrs.Name)" }
$trAnimal = $animalAttrs.TextFrame.TextRange
Set-ParaPrefix $trAnimal 8  "-ADD_PRICE: "         11
Set-ParaPrefix $trAnimal 9  "-SELL_PRICE: "        12
Set-ParaPrefix $trAnimal 10 "-MAINTENANCE_PRICE: " 19
Set-ParaPrefix $trAnimal 11 "-MONEY_GAIN: "        12

# --- Shape "ZoneTexte 327" (Enclosure class attributes) ---
$enclosureAttrs = $s.Shapes.Item(89)
if ($enclosureAttrs.Name -ne "ZoneTexte 327") { throw "unexpected shape 89: $($enclosureAttrs.Name)" }
$trEnclosure = $enclosureAttrs.TextFrame.TextRange
Set-ParaPrefix $trEnclosure 3 "-MAX_ANIMAL: "        12
Set-ParaPrefix $trEnclosure 6 "-ADD_PRICE: "         11
Set-ParaPrefix $trEnclosure 7 "-SELL_PRICE: "        12
Set-ParaPrefix $trEnclosure 8 "-MAINTENANCE_PRICE: " 19

# The longer text re-triggers PowerPoint's autofit-width layout for the
# "Enclosure" attribute box, which also nudges its title box, its
# background rectangle, and the glued connector that anchors to its
# left-centre connection point. Reproduce those exact extents.
$enclosureAttrs.Left = 423.4696197509766
$enclosureAttrs.Width = 81.52984237670898

$enclosureTitle = $s.Shapes.Item(86)
if ($enclosureTitle.Name -ne "ZoneTexte 323") { throw "unexpected shape 86: $($enclosureTitle.Name)" }
$enclosureTitle.Width = 62.84669303894049

$enclosureBody = $s.Shapes.Item(87)
if ($enclosureBody.Name -ne "Rectangle 324") { throw "unexpected shape 87: $($enclosureBody.Name)" }
$enclosureBody.Width = 79.12298965454107

$enclosureConnector = $s.Shapes.Item(82)
if ($enclosureConnector.Name -ne "Connecteur droit 312") { throw "unexpected shape 82: $($enclosureConnector.Name)" }
$enclosureConnector.Width = 45.494802474975586
